$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title "Beschreibung der Kategorien": make it bold (both the paragraph
#    mark run properties and the run itself get <w:b/><w:bCs/>).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

$titleXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Beschreibung der Kategorien</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$titleRange.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2) "Zusaetzliche Funktionen" block -> "Anpassung der Funktionen" block:
#    - "Zusätzliche Funktionen: " / "Erweiterte Funktionen..." / "Bestehende
#      Funktionen verbessern: " / "Verbesserung / Optimierung..." becomes
#    - "Anpassung der Funktionen: " / "Erweiterte Funktionen..." /
#      "Verbesserung / Optimierung..." / "Hinzufügen und Bearbeiten von
#      Funktionen " (new paragraph).
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Zus*tzliche Funktionen*") {
        $startPara = $d.Paragraphs.Item($i)
    }
    if ($txt -like "Verbesserung / Optimierung der bereits bestehenden*") {
        $endPara = $d.Paragraphs.Item($i)
        break
    }
}

$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$blockXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r><w:t>Anpassung der Funktionen</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Erweiterte Funktionen, die über die grundlegenden Funktionen hinausgehen und den Benutzern zusätzlichen Nutzen bieten</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Verbesserung / Optimierung der bereits bestehenden Funktionalitäten</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Hinzufügen und Bearbeiten von Funktionen </w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$blockRange.InsertXML($blockXml)

Write-Host "Edits applied."
